$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
  # Row 92: Whinier than the Sword | Enchanted Koppranickel Ink
  $ws.Range("H92").Value = 443.76923
  $ws.Range("I92").Value = 449.68182
  $ws.Range("J92").Value = 411.25
  $ws.Range("K92").Value = 449.68182
  $ws.Range("L92").Value = 411.25
  $ws.Range("M92").Value = 798.31818
  $ws.Range("N92").Value = -2907.25
  # Row 101: Edge of the Arcane | Cunning Craftsman's Tea
  $ws.Range("H101").Value = 572.7273
  $ws.Range("I101").Value = 339.375
  $ws.Range("J101").Value = 1195
  $ws.Range("K101").Value = 1018.125
  $ws.Range("L101").Value = 3585
  $ws.Range("M101").Value = 603.875
  $ws.Range("N101").Value = -6829
  # Row 103: Let Loose the Juice | Persimmon Tannin
  $ws.Range("H103").Value = 1831.1428
  $ws.Range("J103").Value = 556
  $ws.Range("L103").Value = 1668
  $ws.Range("N103").Value = -2840
  # Row 107: Another Man's Ink | Enchanted Truegold Ink
  $ws.Range("H107").Value = 602.64
  $ws.Range("I107").Value = 607.26086
  $ws.Range("J107").Value = 549.5
  $ws.Range("K107").Value = 607.26086
  $ws.Range("L107").Value = 549.5
  $ws.Range("M107").Value = 1312.73914
  $ws.Range("N107").Value = -4389.5
  # Row 113: Amaro Kart | Starch Glue
  $ws.Range("H113").Value = 113222.78
  $ws.Range("I113").Value = 202201
  $ws.Range("J113").Value = 2000
  $ws.Range("K113").Value = 202201
  $ws.Range("L113").Value = 2000
  $ws.Range("M113").Value = -198947
  $ws.Range("N113").Value = -8508
  # Row 129: Practical Command | Commanding Craftsman's Draught
  $ws.Range("H129").Value = 2864.0417
  $ws.Range("I129").Value = 6395.8237
  $ws.Range("J129").Value = 927.25806
  $ws.Range("K129").Value = 19187.4711
  $ws.Range("L129").Value = 2781.77418
  $ws.Range("M129").Value = -14187.4711
  $ws.Range("N129").Value = -12781.77418
  # Row 132: Fast-forwarding Flora | Growth Formula Lambda
  $ws.Range("H132").Value = 5213460
  $ws.Range("I132").Value = 5957604
  $ws.Range("J132").Value = 4451.8335
  $ws.Range("K132").Value = 17872812
  $ws.Range("L132").Value = 13355.5005
  $ws.Range("M132").Value = -17870282
  $ws.Range("N132").Value = -18415.5005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
  # Row 2: Ain't Got No Ingots | Bronze Ingot
  $ws.Range("H2").Value = 29472.828
  $ws.Range("I2").Value = 923.40625
  $ws.Range("J2").Value = 334000
  $ws.Range("K2").Value = 923.40625
  $ws.Range("L2").Value = 334000
  $ws.Range("M2").Value = -810.40625
  $ws.Range("N2").Value = -334226
  # Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
  $ws.Range("H102").Value = 103053.9
  $ws.Range("I102").Value = 254319.75
  $ws.Range("J102").Value = 2210
  $ws.Range("K102").Value = 254319.75
  $ws.Range("L102").Value = 2210
  $ws.Range("M102").Value = -252697.75
  $ws.Range("N102").Value = -5454
  # Row 116: No Scope | Titanbronze Ingot
  $ws.Range("H116").Value = 29472.828
  $ws.Range("I116").Value = 923.40625
  $ws.Range("J116").Value = 334000
  $ws.Range("K116").Value = 923.40625
  $ws.Range("L116").Value = 334000
  $ws.Range("M116").Value = 1370.59375
  $ws.Range("N116").Value = -338588

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
  # Row 3: Hells Bells | Bronze Ingot
  $ws.Range("H3").Value = 29472.828
  $ws.Range("I3").Value = 923.40625
  $ws.Range("J3").Value = 334000
  $ws.Range("K3").Value = 923.40625
  $ws.Range("L3").Value = 334000
  $ws.Range("M3").Value = -809.40625
  $ws.Range("N3").Value = -334228
  # Row 80: Unbreaker | Titanium Ingot
  $ws.Range("H80").Value = 2904.9312
  $ws.Range("I80").Value = 845.7778
  $ws.Range("J80").Value = 3831.55
  $ws.Range("K80").Value = 845.7778
  $ws.Range("L80").Value = 3831.55
  $ws.Range("M80").Value = 152.2222
  $ws.Range("N80").Value = -5827.55
  # Row 83: Attack on Titanium (L) | Titanium Ingot
  $ws.Range("H83").Value = 2904.9312
  $ws.Range("I83").Value = 845.7778
  $ws.Range("J83").Value = 3831.55
  $ws.Range("K83").Value = 4228.889
  $ws.Range("L83").Value = 19157.75
  $ws.Range("M83").Value = 763.1109999999999
  $ws.Range("N83").Value = -29141.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
  # Row 88: Hold on Adamantite | Adamantite Spear
  $ws.Range("H88").Value = 22671.5
  $ws.Range("J88").Value = 22671.5
  $ws.Range("L88").Value = 22671.5
  $ws.Range("N88").Value = -23483.5
  # Row 91: Spears for Stone Vigilantes (L) | Adamantite Spear
  $ws.Range("H91").Value = 22671.5
  $ws.Range("J91").Value = 22671.5
  $ws.Range("L91").Value = 22671.5
  $ws.Range("N91").Value = -25479.5
  # Row 132: Hull Lotta Damage | Ginseng Lumber
  $ws.Range("H132").Value = 5791.9165
  $ws.Range("I132").Value = 6611.5557
  $ws.Range("J132").Value = 3333
  $ws.Range("K132").Value = 19834.6671
  $ws.Range("L132").Value = 9999
  $ws.Range("M132").Value = -17304.6671
  $ws.Range("N132").Value = -15059

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
  # Row 58: Bread in the Clouds | La Noscean Toast
  $ws.Range("H58").Value = 3400
  $ws.Range("J58").Value = 3400
  $ws.Range("L58").Value = 10200
  $ws.Range("N58").Value = -10456
  # Row 98: Sweet Kiss of Death | Rice Vinegar
  $ws.Range("H98").Value = 72380.92999999999
  $ws.Range("J98").Value = 84354.414
  $ws.Range("L98").Value = 253063.242
  $ws.Range("N98").Value = -256059.242

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
  # Row 113: Copious Crystal Cannons | Manasilver Nugget
  $ws.Range("H113").Value = 1628.6666
  $ws.Range("I113").Value = 1357.5
  $ws.Range("J113").Value = 1899.8334
  $ws.Range("K113").Value = 1357.5
  $ws.Range("L113").Value = 1899.8334
  $ws.Range("M113").Value = 812.5
  $ws.Range("N113").Value = -6239.8334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
  # Row 7: Tan Before the Ban | Leather
  $ws.Range("H7").Value = 6426
  $ws.Range("I7").Value = 2972.6667
  $ws.Range("K7").Value = 2972.6667
  $ws.Range("M7").Value = -2860.6667
  # Row 40: Best Served Toad | Toad Leather
  $ws.Range("H40").Value = 23010.256
  $ws.Range("I40").Value = 35767.45
  $ws.Range("K40").Value = 35767.45
  $ws.Range("M40").Value = -35631.45
  # Row 55: It's Not a Job, It's a Calling | Peiste Leather
  $ws.Range("H55").Value = 975.5714
  $ws.Range("J55").Value = 673.5
  $ws.Range("L55").Value = 673.5
  $ws.Range("N55").Value = -1019.5
  # Row 61: Spelling Me Softly | Raptor Leather
  $ws.Range("H61").Value = 1791.625
  $ws.Range("I61").Value = 1661.1
  $ws.Range("J61").Value = 2009.1666
  $ws.Range("K61").Value = 1661.1
  $ws.Range("L61").Value = 2009.1666
  $ws.Range("M61").Value = -1459.1
  $ws.Range("N61").Value = -2413.1666
  # Row 113: Peace in Rest | Atrociraptor Leather
  $ws.Range("H113").Value = 1791.625
  $ws.Range("I113").Value = 1661.1
  $ws.Range("J113").Value = 2009.1666
  $ws.Range("K113").Value = 1661.1
  $ws.Range("L113").Value = 2009.1666
  $ws.Range("M113").Value = 508.9000000000001
  $ws.Range("N113").Value = -6349.1666
  # Row 126: Battered Books | Saiga Leather
  $ws.Range("H126").Value = 6426
  $ws.Range("I126").Value = 2972.6667
  $ws.Range("K126").Value = 8918.000100000001
  $ws.Range("M126").Value = -6448.000100000001
  # Row 132: Tenets of Tanning | Silver Lobo Leather
  $ws.Range("H132").Value = 5889.5884
  $ws.Range("I132").Value = 6515.5557
  $ws.Range("J132").Value = 5185.375
  $ws.Range("K132").Value = 19546.6671
  $ws.Range("L132").Value = 15556.125
  $ws.Range("M132").Value = -17016.6671
  $ws.Range("N132").Value = -20616.125

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
  # Row 33: I'll Be Your Wailer Today | Velveteen Wedge Cap of Gathering
  $ws.Range("H33").Value = 59953.668
  $ws.Range("J33").Value = 85021
  $ws.Range("L33").Value = 85021
  $ws.Range("N33").Value = -85521
  # Row 34: He's Got Legs | Velveteen Sarouel
  $ws.Range("H34").Value = 16200
  $ws.Range("J34").Value = 16200
  $ws.Range("L34").Value = 16200
  $ws.Range("N34").Value = -16606
  # Row 36: Put a Lid on It | Velveteen Wedge Cap of Gathering
  $ws.Range("H36").Value = 59953.668
  $ws.Range("J36").Value = 85021
  $ws.Range("L36").Value = 85021
  $ws.Range("N36").Value = -85521
  # Row 37: Bet You Anything | Velveteen Sarouel of Gathering
  $ws.Range("H37").Value = 13228.286
  $ws.Range("J37").Value = 13799.667
  $ws.Range("L37").Value = 13799.667
  $ws.Range("N37").Value = -14205.667
  # Row 43: Walk Softly and Carry a Big Halberd | Velveteen Dress Shoes
  $ws.Range("H43").Value = 16200
  $ws.Range("J43").Value = 16200
  $ws.Range("L43").Value = 16200
  $ws.Range("N43").Value = -16498
  # Row 56: Full Moon Fever | Felt Chausses
  $ws.Range("H56").Value = 15163.333
  $ws.Range("I56").Value = 3750
  $ws.Range("J56").Value = 37990
  $ws.Range("K56").Value = 3750
  $ws.Range("L56").Value = 37990
  $ws.Range("M56").Value = -3036
  $ws.Range("N56").Value = -39418
  # Row 136: Weaving the Envelope | Sarcenet Cloth
  $ws.Range("H136").Value = 1485.625
  $ws.Range("I136").Value = 505.08
  $ws.Range("J136").Value = 2276.3872
  $ws.Range("K136").Value = 1515.24
  $ws.Range("L136").Value = 6829.1616
  $ws.Range("M136").Value = 1034.76
  $ws.Range("N136").Value = -11929.1616
